$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old range entirely (A1:E11) before writing the new, smaller table
$ws.Range("A1:E11").Clear()

# Write new headers
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"

# Write new data row
$ws.Range("A2").Value = 81
$ws.Range("B2").Value = 40
$ws.Range("C2").Value = 12987
$ws.Range("D2").Value = 0.2889342308044434
